$wb = $excel.ActiveWorkbook

# The workbook currently has two tabs, in this physical order:
#   slot 1) "hotel_info"  (rId1 / sheet1.xml) - header row + one data row
#   slot 2) "review_info" (rId2 / sheet2.xml) - header row only
#
# The target edit keeps the two underlying sheet slots (and therefore their
# r:id / sheetId pairing) exactly where they are, but swaps which logical
# table lives in which slot:
#   slot 1 (rId1) becomes "review_info" (header row only, no data rows)
#   slot 2 (rId2) becomes "hotel_info"  (header row + one data row, with a
#                                        new "State" column inserted right
#                                        after "Hotel_Name", valued
#                                        "Louisiana")

$sheetSlot1 = $wb.Worksheets.Item(1)
$sheetSlot2 = $wb.Worksheets.Item(2)

# --- Move the original hotel_info data row from slot 1 into its new home
#     in slot 2 *before* slot 1 gets cleared/repurposed. Using Copy/
#     PasteSpecial (instead of re-typing the values) preserves each cell's
#     original data type (number vs text) exactly, including the
#     numeric-looking text values ("91", "1", "92").
#
# Original slot-1 columns: STR, Hotel_Name, City, Zip, TA_ReviewURL,
#   Tripadvisor_Hotel_Name, English_Reviews_num, Local_Rank, Total_Reviews_num
# New slot-2 columns:      STR, Hotel_Name, State, City, Zip, TA_ReviewURL,
#   Tripadvisor_Hotel_Name, English_Reviews_num, Local_Rank, Total_Reviews_num
#
# So columns A:B carry straight across, a new "State" column is inserted at
# C, and the remaining original columns C:I shift right to D:J.
$sheetSlot1.Range("A2:B2").Copy()
$sheetSlot2.Range("A2:B2").PasteSpecial()

$sheetSlot1.Range("C2:I2").Copy()
$sheetSlot2.Range("D2:J2").PasteSpecial()

$sheetSlot2.Cells.Item(2, 3).Value = "Louisiana"

# Rename both sheets through temporary names first so the swap doesn't
# collide with the sheet that currently holds the target name.
$sheetSlot1.Name = "__tmp_slot1__"
$sheetSlot2.Name = "__tmp_slot2__"

# --- Turn slot 1 into the new "review_info" sheet -------------------------
$sheetSlot1.Cells.Clear()
$sheetSlot1.Name = "review_info"

$reviewHeaders = @(
    "STR", "reviewer_ID", "reviewer_name", "Review_ID", "Date_of_scraping",
    "ReviewURL", "Tripadvisor_gcode", "Tripadvisor_dcode", "Tripadvisor_rcode",
    "review_date", "review_title", "review_content", "review_rating",
    "trip_month", "trip_purpose", "value", "rooms", "Location",
    "Cleanliness", "Sleep Quality", "Service", "Picture(yes=1)",
    "respondent", "response_date", "response_text"
)
for ($i = 0; $i -lt $reviewHeaders.Length; $i++) {
    $sheetSlot1.Cells.Item(1, $i + 1).Value = $reviewHeaders[$i]
}

# --- Turn slot 2 into the new "hotel_info" sheet ---------------------------
# Clear out the leftover review_info header cells beyond the columns the
# hotel_info table uses (K1:Y1), but keep the data row (A2:J2) we just
# pasted in above.
$sheetSlot2.Range("K1:Y1").Clear()
$sheetSlot2.Name = "hotel_info"

$hotelHeaders = @(
    "STR", "Hotel_Name", "State", "City", "Zip", "TA_ReviewURL",
    "Tripadvisor_Hotel_Name", "English_Reviews_num", "Local_Rank",
    "Total_Reviews_num"
)
for ($i = 0; $i -lt $hotelHeaders.Length; $i++) {
    $sheetSlot2.Cells.Item(1, $i + 1).Value = $hotelHeaders[$i]
}
